# Update the "Förändrad" (Changed) date column (C) for rows 2-37
# from 45648 (2024-12-22) to 45649 (2024-12-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 37; $row++) {
    $ws.Cells.Item($row, 3).Value = 45649
}
